# Auto-generated edit script applying market-price refresh to Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1381.3334
$ws.Range("I43").Value = 1572
$ws.Range("K43").Value = 1572
$ws.Range("M43").Value = -1503
$ws.Range("H53").Value = 215
$ws.Range("I53").Value = 77.666664
$ws.Range("K53").Value = 77.666664
$ws.Range("M53").Value = 559.333336
$ws.Range("H74").Value = 10709.857
$ws.Range("I74").Value = 10709.857
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 10709.857
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 10709.857
$ws.Range("I77").Value = 10709.857
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 53549.285
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 4016.818
$ws.Range("I86").Value = 4666.6665
$ws.Range("K86").Value = 4666.6665
$ws.Range("M86").Value = -3543.6665
$ws.Range("H89").Value = 4016.818
$ws.Range("I89").Value = 4666.6665
$ws.Range("K89").Value = 23333.3325
$ws.Range("M89").Value = -17717.3325
$ws.Range("H98").Value = 4515.5
$ws.Range("I98").Value = 3527.8333
$ws.Range("J98").Value = 5997
$ws.Range("K98").Value = 3527.8333
$ws.Range("L98").Value = 5997
$ws.Range("M98").Value = -2029.8333
$ws.Range("N98").Value = -8993
$ws.Range("H122").Value = 4515.5
$ws.Range("I122").Value = 3527.8333
$ws.Range("J122").Value = 5997
$ws.Range("K122").Value = 10583.4999
$ws.Range("L122").Value = 17991
$ws.Range("M122").Value = -8133.499899999999
$ws.Range("N122").Value = -22891
$ws.Range("H135").Value = 575.6
$ws.Range("I135").Value = 391.29166
$ws.Range("J135").Value = 4999
$ws.Range("K135").Value = 3521.62494
$ws.Range("L135").Value = 44991
$ws.Range("M135").Value = -986.6249399999997
$ws.Range("N135").Value = -50061
$ws.Range("H137").Value = 2341.6
$ws.Range("I137").Value = 2177
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 6531
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -3981
$ws.Range("N137").Value = -14100
$ws.Range("H138").Value = 4409.879
$ws.Range("J138").Value = 4572.393
$ws.Range("L138").Value = 13717.179
$ws.Range("N138").Value = -23997.179
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9435.588
$ws.Range("I32").Value = 8921.226000000001
$ws.Range("K32").Value = 8921.226000000001
$ws.Range("M32").Value = -8634.226000000001
$ws.Range("H61").Value = 2753
$ws.Range("I61").Value = 2748.1667
$ws.Range("K61").Value = 2748.1667
$ws.Range("M61").Value = -2536.1667
$ws.Range("H122").Value = 1703.8572
$ws.Range("I122").Value = 1450.3077
$ws.Range("K122").Value = 4350.9231
$ws.Range("M122").Value = -1900.9231
$ws.Range("H136").Value = 2753
$ws.Range("I136").Value = 2748.1667
$ws.Range("K136").Value = 8244.500100000001
$ws.Range("M136").Value = -5694.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1788.5
$ws.Range("I134").Value = 1788.6666
$ws.Range("J134").Value = 1788
$ws.Range("K134").Value = 5365.9998
$ws.Range("L134").Value = 5364
$ws.Range("M134").Value = -2830.9998
$ws.Range("N134").Value = -10434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2052.2856
$ws.Range("I31").Value = 1248.3636
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1248.3636
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -953.3635999999999
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 2052.2856
$ws.Range("I34").Value = 1248.3636
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1248.3636
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -1046.3636
$ws.Range("N34").Value = -5404
$ws.Range("H86").Value = 5956.722
$ws.Range("I86").Value = 4856.625
$ws.Range("K86").Value = 4856.625
$ws.Range("M86").Value = -3733.625
$ws.Range("H89").Value = 5956.722
$ws.Range("I89").Value = 4856.625
$ws.Range("K89").Value = 24283.125
$ws.Range("M89").Value = -18667.125
$ws.Range("H107").Value = 2119.75
$ws.Range("I107").Value = 1715.4445
$ws.Range("K107").Value = 1715.4445
$ws.Range("M107").Value = 204.5554999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 44.272728
$ws.Range("J33").Value = 13.5
$ws.Range("L33").Value = 81
$ws.Range("N33").Value = -647
$ws.Range("H104").Value = 4946.0557
$ws.Range("J104").Value = 4946.0557
$ws.Range("L104").Value = 14838.1671
$ws.Range("N104").Value = -20080.1671
$ws.Range("H116").Value = 4999.6816
$ws.Range("I116").Value = 4994
$ws.Range("K116").Value = 14982
$ws.Range("M116").Value = -11540
$ws.Range("H124").Value = 4904.5
$ws.Range("I124").Value = 2900
$ws.Range("J124").Value = 4999.952
$ws.Range("K124").Value = 8700
$ws.Range("L124").Value = 14999.856
$ws.Range("M124").Value = -3790
$ws.Range("N124").Value = -24819.856
$ws.Range("H136").Value = 2320.7646
$ws.Range("I136").Value = 1988.25
$ws.Range("J136").Value = 2423.077
$ws.Range("K136").Value = 5964.75
$ws.Range("L136").Value = 7269.231000000001
$ws.Range("M136").Value = -864.75
$ws.Range("N136").Value = -17469.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8500
$ws.Range("I70").Value = 8500
$ws.Range("K70").Value = 8500
$ws.Range("M70").Value = -8230
$ws.Range("H73").Value = 8500
$ws.Range("I73").Value = 8500
$ws.Range("K73").Value = 8500
$ws.Range("M73").Value = -7564
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 4874
$ws.Range("I102").Value = 4874
$ws.Range("K102").Value = 4874
$ws.Range("M102").Value = -3252
$ws.Range("H113").Value = 2083.7368
$ws.Range("I113").Value = 1548
$ws.Range("K113").Value = 1548
$ws.Range("M113").Value = 622
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7499.1816
$ws.Range("I22").Value = 6730.3335
$ws.Range("K22").Value = 6730.3335
$ws.Range("M22").Value = -6435.3335
$ws.Range("H27").Value = 7499.1816
$ws.Range("I27").Value = 6730.3335
$ws.Range("K27").Value = 6730.3335
$ws.Range("M27").Value = -6623.3335
$ws.Range("H136").Value = 3229.8572
$ws.Range("I136").Value = 2175.875
$ws.Range("K136").Value = 6527.625
$ws.Range("M136").Value = -3977.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 28888.75
$ws.Range("J54").Value = 33518.332
$ws.Range("L54").Value = 33518.332
$ws.Range("N54").Value = -34558.332
$ws.Range("H81").Value = 2412.9167
$ws.Range("I81").Value = 1814.091
$ws.Range("K81").Value = 3628.182
$ws.Range("M81").Value = -2567.182
$ws.Range("H84").Value = 2412.9167
$ws.Range("I84").Value = 1814.091
$ws.Range("K84").Value = 18140.91
$ws.Range("M84").Value = -12836.91
$ws.Range("H96").Value = 2200
$ws.Range("I96").Value = 2250
$ws.Range("J96").Value = 2150
$ws.Range("K96").Value = 2250
$ws.Range("L96").Value = 2150
$ws.Range("M96").Value = -877
$ws.Range("N96").Value = -4896
$ws.Range("H113").Value = 713.4
$ws.Range("I113").Value = 740.375
$ws.Range("J113").Value = 695.4167
$ws.Range("K113").Value = 2221.125
$ws.Range("L113").Value = 2086.2501
$ws.Range("M113").Value = -51.125
$ws.Range("N113").Value = -6426.2501
$ws.Range("H122").Value = 1113.2727
$ws.Range("I122").Value = 1124.5
$ws.Range("K122").Value = 3373.5
$ws.Range("M122").Value = -923.5
$ws.Range("H126").Value = 1690.762
$ws.Range("I126").Value = 1635.3334
$ws.Range("K126").Value = 4906.0002
$ws.Range("M126").Value = -2436.0002
$ws.Range("H132").Value = 3645.4375
$ws.Range("I132").Value = 3206.6667
$ws.Range("K132").Value = 9620.000100000001
$ws.Range("M132").Value = -7090.000100000001
